$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 / J column: "nghithucnha" label replaced by a numeric code
$ws.Range("J2").Value = 11110000

# I3: "Nộp Đơn" -> "Nộp tại đây"
$ws.Range("I3").Value = "Nộp tại đây"

# F5: remove the stray "Bán tài khoản và phần mềm uy tín" text
$ws.Range("F5").ClearContents()

# I7: "Đang Cập Nhập" -> "Xem ảnh"
$ws.Range("I7").Value = "Xem ảnh"

# H7: add the new cover image filename
$ws.Range("H7").Value = "bonmang2026.jpg"

# G7: add the new Google Drive share link (Bổn Mạng Nhà 2026 album)
$ws.Range("G7").Value = "https://drive.google.com/drive/folders/15vvwM8k_XceMiOjOIIpUWbvNwV6kRFwR?usp=sharing"

# G8: append "?usp=sharing" to the existing Tu Lieu 2025 drive link
$ws.Range("G8").Value = "https://drive.google.com/drive/folders/1IVawCMt9xO_6Cnvzh2S28Q6U66pePz7e?usp=sharing"

# New column J width (col 10), stored width 9.5
$ws.Columns.Item(10).ColumnWidth = 8.666666666666666

# Update the saved selection to G8
$ws.Range("G8").Select()
